# Apply logbook entries 8 and 9 to the Documentation Logbook worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string text for the two new entries.
$title8 = "Finished implementing player movement. Began implementing border generation and movement."
$desc8 = "Fixed player movement bug where the character wouldn" + [char]0x2019 + "t always move in the desired direction. Began to implement Border() class and attempted to make them visible."

# Row 12: entry No. 8 dated 24/01/2021 (serial 44220) with Title/Description text.
$ws.Range("B12").Value = 8
$ws.Range("C12").Value = [DateTime]::FromOADate(44220)
$ws.Range("D12").Value = $title8
$ws.Range("E12").Value = $desc8

# Row 13: only the No. column is filled in for entry 9 so far.
$ws.Range("B13").Value = 9

# The wrapped text in row 12 now spans two lines, so the row grows to fit it.
$ws.Rows.Item(12).RowHeight = 43.5

# Leave the selection on B13, matching the saved state in the workbook.
$ws.Range("B13").Select() | Out-Null
